$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped on Fri Oct 20 16:08:11 UTC 2023
# (includes two coin-row swaps: Stellar/EthereumClassic and PaxDollar/BitcoinSV)
$ws.Range('D2').Value = '29.474.00'
$ws.Range('E2').Value = '  +2.76%  '

$ws.Range('D3').Value = '1.604.04'
$ws.Range('E3').Value = '  +2.66%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').Value = '''212.74'

$ws.Range('D6').Value = '''0.523'
$ws.Range('E6').Value = '  +7.49%  '

$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('D8').Value = '''26.87'
$ws.Range('E8').Value = '  +8.00%  '

$ws.Range('D9').Value = '''43.53'
$ws.Range('E9').Value = '  -1.18%  '

$ws.Range('D10').Value = '''0.251'
$ws.Range('E10').Value = '  +2.19%  '

$ws.Range('D11').Value = '''0.0599'
$ws.Range('E11').Value = '  +2.48%  '

$ws.Range('E12').Value = '  +1.92%  '

$ws.Range('D13').Value = '1.832.44'
$ws.Range('E13').Value = '  +2.68%  '

$ws.Range('D14').Value = '1.595.21'
$ws.Range('E14').Value = '  +1.92%  '

$ws.Range('D15').Value = '29.513.08'
$ws.Range('E15').Value = '  +2.99%  '

$ws.Range('E16').Value = '  +3.97%  '

$ws.Range('E17').Value = '  +1.96%  '

$ws.Range('E18').Value = '  +3.14%  '

$ws.Range('D19').Value = '''243.00'
$ws.Range('E19').Value = '  +5.57%  '

$ws.Range('D20').Value = '''7.63'
$ws.Range('E20').Value = '  +3.51%  '

$ws.Range('D21').Value = '0.0₃0691'
$ws.Range('E21').Value = '  +2.57%  '

$ws.Range('E22').Value = '  -0.03%  '

$ws.Range('D23').Value = '''3.99'
$ws.Range('E23').Value = '  +1.65%  '

$ws.Range('D24').Value = '''9.16'
$ws.Range('E24').Value = '  +1.73%  '

$ws.Range('E25').Value = '  +0.32%  '

$ws.Range('D26').Value = '''154.63'
$ws.Range('E26').Value = '  +2.49%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''15.32'
$ws.Range('E27').Value = '  +3.59%  '

$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '''0.109'
$ws.Range('E28').Value = '  +5.61%  '

$ws.Range('D29').Value = '''6.39'
$ws.Range('E29').Value = '  +2.59%  '

$ws.Range('E30').Value = '  +0.08%  '

$ws.Range('E32').Value = '  -0.29%  '

$ws.Range('E33').Value = '  +1.80%  '

$ws.Range('D34').Value = '1.418.47'

$ws.Range('D35').Value = '''3.10'
$ws.Range('E35').Value = '  +3.86%  '

$ws.Range('E36').Value = '  -1.78%  '

$ws.Range('E37').Value = '  +2.69%  '

$ws.Range('D38').Value = '''2.80'
$ws.Range('E38').Value = '  +5.38%  '

$ws.Range('E39').Value = '  +0.29%  '

$ws.Range('E40').Value = '  +2.25%  '

$ws.Range('E41').Value = '  +3.42%  '

$ws.Range('E42').Value = '  +0.37%  '

$ws.Range('B43').Value = 'BitcoinSV'
$ws.Range('C43').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D43').Value = '''53.87'
$ws.Range('E43').Value = '  +23.88%  '

$ws.Range('D44').Value = '''0.0484'
$ws.Range('E44').Value = '  +4.68%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '''0.998'
$ws.Range('E45').Value = '  +0.02%  '

$ws.Range('D46').Value = '''0.793'
$ws.Range('E46').Value = '  +2.22%  '

$ws.Range('D47').Value = '''65.57'
$ws.Range('E47').Value = '  +2.54%  '

$ws.Range('D48').Value = '''5.29'
$ws.Range('E48').Value = '  +0.18%  '

$ws.Range('D49').Value = '1.743.02'
$ws.Range('E49').Value = '  +2.84%  '

$ws.Range('D50').Value = '''86.52'
$ws.Range('E50').Value = '  +1.43%  '

$ws.Range('D51').Value = '''0.831'
$ws.Range('E51').Value = '  -4.38%  '
